# "Responsive para pantallas pequeñas y retoques"
# - Add three new task rows (13, 14, 15) to the tracker sheet.
# - Tighten several row heights (content made more compact / responsive).
# - Move the active selection to F19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 17-19 (tasks #13, #14, #15) -----------------------------
$ws.Range("A17").Value = 13
$ws.Range("B17").Value = "Modificacion final al header y adicion de Logo"
$ws.Range("C17").Value = "Samuel"
$ws.Range("D17").Value = 45753
$ws.Range("E17").Value = 45753
$ws.Range("F17").Value = "✅ Hecho"

$ws.Range("A18").Value = 14
$ws.Range("B18").Value = "Responsive y ajustes"
$ws.Range("C18").Value = "Jon"
$ws.Range("D18").Value = 45753
$ws.Range("E18").Value = 45753
$ws.Range("F18").Value = "✅ Hecho"

$ws.Range("A19").Value = 15
$ws.Range("B19").Value = "Entidad Relacion"
$ws.Range("C19").Value = "Samuel"
$ws.Range("D19").Value = 45753
$ws.Range("F19").Value = "🕓 Pendiente"

# --- Responsive row-height retouches -----------------------------------
$ws.Rows(2).RowHeight = 30.75
$ws.Rows(4).RowHeight = 30
$ws.Rows(5).RowHeight = 30
$ws.Rows(8).AutoFit()
$ws.Rows(11).RowHeight = 60
$ws.Rows(12).AutoFit()
$ws.Rows(13).AutoFit()
$ws.Rows(14).RowHeight = 150
$ws.Rows(17).RowHeight = 30
$ws.Rows(24).RowHeight = 45

# --- Selection moves to the newest edited cell --------------------------
$ws.Range("F19").Select() | Out-Null
